# "cleaned defensive actions data"
#
# The sheet was produced by pandas.to_excel() with a 2-row multi-index
# header: row 1 held the raw (mostly "Unnamed: N_level_0") top level and
# row 2 held the real column names, with a handful of top-level groups
# (Tackles / Challenges / Blocks) merged across their sub-columns. This
# edit "flattens" that back into a single, real header on row 1, hides
# the now-redundant old header row (2) and the blank spacer row (3),
# un-merges the old merged header cells, hides the trailing totals row
# (20), and fills in a few stray blank numeric cells in column O so the
# data range is fully rectangular.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the old merged header groupings first -- a merged cell's
#        value only "sticks" on its top-left member, so the individual
#        sub-headers beneath H1/M1/Q1 can't be written until the ranges
#        are split back into ordinary cells. ----------------------------
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- 2. Replace row 1 with the real, flattened column headers ----------
$headers = @{
    "A1" = "Player ID"
    "B1" = "Player"
    "C1" = "#"
    "D1" = "Nation"
    "E1" = "Pos"
    "F1" = "Age"
    "G1" = "90s"
    "H1" = "Tkl"
    "I1" = "TklW"
    "J1" = "Def 3rd"
    "K1" = "Mid 3rd"
    "L1" = "Att 3rd"
    "M1" = "Cha"
    "N1" = "Att"
    "O1" = "Tkl%"
    "P1" = "Lost"
    "Q1" = "Blocks"
    "R1" = "Sh"
    "S1" = "Pass"
    "T1" = "Int"
    "U1" = "Tkl+Int"
    "V1" = "Clr"
    "W1" = "Err"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- 3. Hide the stale duplicate header row, the blank spacer row, and
#        the trailing totals row -----------------------------------------
$ws.Range("A2").EntireRow.Hidden = $true
$ws.Range("A3").EntireRow.Hidden = $true
$ws.Range("A20").EntireRow.Hidden = $true

# --- 4. Fill in the missing blank numeric cells in column O -------------
$oCells = @(6, 8, 10, 11, 13, 15, 16, 17, 19)
foreach ($r in $oCells) {
    $ws.Range("O$r").Value = 0
}

# --- 5. Restore the saved selection --------------------------------------
[void]$ws.Range("O21").Select()
